# Scheduled market-data refresh: update currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H:N) for the affected leve rows across sheets, per the upstream run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2314.4
$ws.Range("I28").Value = 2369.9
$ws.Range("J28").Value = 2203.4
$ws.Range("K28").Value = 2369.9
$ws.Range("L28").Value = 2203.4
$ws.Range("M28").Value = -1884.9
$ws.Range("N28").Value = -3173.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 10005280
$ws.Range("I32").Value = 1500
$ws.Range("J32").Value = 11116811
$ws.Range("K32").Value = 1500
$ws.Range("L32").Value = 11116811
$ws.Range("M32").Value = -1174
$ws.Range("N32").Value = -11117463

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 8562
$ws.Range("I74").Value = 6577.25
$ws.Range("J74").Value = 10149.8
$ws.Range("K74").Value = 6577.25
$ws.Range("L74").Value = 10149.8
$ws.Range("M74").Value = -5641.25
$ws.Range("N74").Value = -12021.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 8562
$ws.Range("I77").Value = 6577.25
$ws.Range("J77").Value = 10149.8
$ws.Range("K77").Value = 32886.25
$ws.Range("L77").Value = 50749
$ws.Range("M77").Value = -28206.25
$ws.Range("N77").Value = -60109

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1137.8096
$ws.Range("I92").Value = 1052
$ws.Range("J92").Value = 1309.4286
$ws.Range("K92").Value = 1052
$ws.Range("L92").Value = 1309.4286
$ws.Range("M92").Value = 196
$ws.Range("N92").Value = -3805.4286

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 5853.32
$ws.Range("I100").Value = 4390.4614
$ws.Range("K100").Value = 4390.4614
$ws.Range("M100").Value = -3849.4614

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 1644.8572
$ws.Range("I101").Value = 259.8
$ws.Range("K101").Value = 779.4000000000001
$ws.Range("M101").Value = 842.5999999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 207.58333
$ws.Range("I107").Value = 207.58333
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 207.58333
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1712.41667
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3389.9583
$ws.Range("I141").Value = 3198.158
$ws.Range("K141").Value = 9594.474
$ws.Range("M141").Value = -4414.474

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1004.92
$ws.Range("J88").Value = 967.25
$ws.Range("L88").Value = 967.25
$ws.Range("N88").Value = -1779.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 1004.92
$ws.Range("J91").Value = 967.25
$ws.Range("L91").Value = 967.25
$ws.Range("N91").Value = -3775.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1672.1875
$ws.Range("I110").Value = 1583.6666
$ws.Range("K110").Value = 1583.6666
$ws.Range("M110").Value = 461.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3019.3333
$ws.Range("I132").Value = 2267.6
$ws.Range("K132").Value = 6802.799999999999
$ws.Range("M132").Value = -4272.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2785.7144
$ws.Range("I94").Value = 2617.9092
$ws.Range("K94").Value = 2617.9092
$ws.Range("M94").Value = -2166.9092

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3664.818
$ws.Range("I132").Value = 3664.818
$ws.Range("K132").Value = 10994.454
$ws.Range("M132").Value = -8464.454000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 274.83334
$ws.Range("J23").Value = 287.25
$ws.Range("L23").Value = 861.75
$ws.Range("N23").Value = -1331.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 55.9
$ws.Range("I40").Value = 54.42857
$ws.Range("K40").Value = 217.71428
$ws.Range("M40").Value = -148.71428

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 186566.67
$ws.Range("I80").Value = 223721.2
$ws.Range("J80").Value = 794
$ws.Range("K80").Value = 223721.2
$ws.Range("L80").Value = 794
$ws.Range("M80").Value = -222723.2
$ws.Range("N80").Value = -2790

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 186566.67
$ws.Range("I83").Value = 223721.2
$ws.Range("J83").Value = 794
$ws.Range("K83").Value = 1118606
$ws.Range("L83").Value = 3970
$ws.Range("M83").Value = -1113614
$ws.Range("N83").Value = -13954

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 562.5
$ws.Range("I97").Value = 550
$ws.Range("J97").Value = 575
$ws.Range("K97").Value = 550
$ws.Range("L97").Value = 575
$ws.Range("M97").Value = -54
$ws.Range("N97").Value = -1567

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 6610.5
$ws.Range("I113").Value = 3413.3635
$ws.Range("J113").Value = 18333.334
$ws.Range("K113").Value = 3413.3635
$ws.Range("L113").Value = 18333.334
$ws.Range("M113").Value = -1243.3635
$ws.Range("N113").Value = -22673.334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2799.9644
$ws.Range("I132").Value = 1961.7646
$ws.Range("J132").Value = 4095.3635
$ws.Range("K132").Value = 5885.293799999999
$ws.Range("L132").Value = 12286.0905
$ws.Range("M132").Value = -3355.293799999999
$ws.Range("N132").Value = -17346.0905

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8282.5
$ws.Range("I7").Value = 12300
$ws.Range("K7").Value = 12300
$ws.Range("M7").Value = -12188

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1128.4
$ws.Range("I22").Value = 1048.3334
$ws.Range("J22").Value = 1248.5
$ws.Range("K22").Value = 1048.3334
$ws.Range("L22").Value = 1248.5
$ws.Range("M22").Value = -753.3334
$ws.Range("N22").Value = -1838.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1128.4
$ws.Range("I27").Value = 1048.3334
$ws.Range("J27").Value = 1248.5
$ws.Range("K27").Value = 1048.3334
$ws.Range("L27").Value = 1248.5
$ws.Range("M27").Value = -941.3334
$ws.Range("N27").Value = -1462.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1197
$ws.Range("I46").Value = 407.5
$ws.Range("J46").Value = 1535.3572
$ws.Range("K46").Value = 407.5
$ws.Range("L46").Value = 1535.3572
$ws.Range("M46").Value = -219.5
$ws.Range("N46").Value = -1911.3572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H118").Value = 48922.69
$ws.Range("J118").Value = 48922.69
$ws.Range("L118").Value = 48922.69
$ws.Range("N118").Value = -52236.69

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 8282.5
$ws.Range("I126").Value = 12300
$ws.Range("K126").Value = 36900
$ws.Range("M126").Value = -34430

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8482.866
$ws.Range("I81").Value = 21729.2
$ws.Range("J81").Value = 1859.7
$ws.Range("K81").Value = 43458.4
$ws.Range("L81").Value = 3719.4
$ws.Range("M81").Value = -42397.4
$ws.Range("N81").Value = -5841.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 8482.866
$ws.Range("I84").Value = 21729.2
$ws.Range("J84").Value = 1859.7
$ws.Range("K84").Value = 217292
$ws.Range("L84").Value = 18597
$ws.Range("M84").Value = -211988
$ws.Range("N84").Value = -29205

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 893.41174
$ws.Range("I100").Value = 798.75
$ws.Range("K100").Value = 1597.5
$ws.Range("M100").Value = -1056.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 560.6667
$ws.Range("I107").Value = 593.375
$ws.Range("J107").Value = 299
$ws.Range("K107").Value = 1780.125
$ws.Range("L107").Value = 897
$ws.Range("M107").Value = 139.875
$ws.Range("N107").Value = -4737

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2530.457
$ws.Range("I132").Value = 1828.5714
$ws.Range("K132").Value = 5485.7142
$ws.Range("M132").Value = -2955.7142
